$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A23").Value = "0036.020320/2025-35"
$ws.Range("B23").Value = "Congresso do Conselho Nacional de Secretario em Brasília"
